$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_9_2_24"
$ws.Cells.Item(2, 2).Value = 0.3695088074892304
$ws.Cells.Item(2, 3).Value = -0.1106217635489666
$ws.Cells.Item(2, 4).Value = -0.9762232811601388
$ws.Cells.Item(2, 5).Value = -0.578781877285212
$ws.Cells.Item(2, 6).Value = 0.6977679133415222
$ws.Cells.Item(2, 7).Value = 1.126189827919006
$ws.Cells.Item(2, 8).Value = 2.817254781723022
$ws.Cells.Item(2, 9).Value = 1.921980738639832

$ws.Cells.Item(3, 1).Value = "model_9_2_23"
$ws.Cells.Item(3, 2).Value = 0.3790565085579682
$ws.Cells.Item(3, 3).Value = -0.1001175054870564
$ws.Cells.Item(3, 4).Value = -0.9266719344703946
$ws.Cells.Item(3, 5).Value = -0.5468473346388312
$ws.Cells.Item(3, 6).Value = 0.6872014999389648
$ws.Cells.Item(3, 7).Value = 1.115538358688354
$ws.Cells.Item(3, 8).Value = 2.746615409851074
$ws.Cells.Item(3, 9).Value = 1.883104205131531

$ws.Cells.Item(4, 1).Value = "model_9_2_22"
$ws.Cells.Item(4, 2).Value = 0.3799867995282954
$ws.Cells.Item(4, 3).Value = -0.105999169459599
$ws.Cells.Item(4, 4).Value = -0.915921015579485
$ws.Cells.Item(4, 5).Value = -0.5435168624632678
$ws.Cells.Item(4, 6).Value = 0.6861719489097595
$ws.Cells.Item(4, 7).Value = 1.12150239944458
$ws.Cells.Item(4, 8).Value = 2.731289386749268
$ws.Cells.Item(4, 9).Value = 1.879049777984619

$ws.Cells.Item(5, 1).Value = "model_9_2_21"
$ws.Cells.Item(5, 2).Value = 0.3893235075789721
$ws.Cells.Item(5, 3).Value = -0.09959722668754245
$ws.Cells.Item(5, 4).Value = -0.8597958380747621
$ws.Cells.Item(5, 5).Value = -0.5097650533934832
$ws.Cells.Item(5, 6).Value = 0.6758389472961426
$ws.Cells.Item(5, 7).Value = 1.115010738372803
$ws.Cells.Item(5, 8).Value = 2.651278495788574
$ws.Cells.Item(5, 9).Value = 1.837960720062256

$ws.Cells.Item(6, 1).Value = "model_9_2_20"
$ws.Cells.Item(6, 2).Value = 0.3896898229213331
$ws.Cells.Item(6, 3).Value = -0.107390708858921
$ws.Cells.Item(6, 4).Value = -0.8509379602131286
$ws.Cells.Item(6, 5).Value = -0.5083196856116836
$ws.Cells.Item(6, 6).Value = 0.6754335761070251
$ws.Cells.Item(6, 7).Value = 1.122913360595703
$ws.Cells.Item(6, 8).Value = 2.638651132583618
$ws.Cells.Item(6, 9).Value = 1.836201310157776

$ws.Cells.Item(7, 1).Value = "model_9_2_19"
$ws.Cells.Item(7, 2).Value = 0.394599369548678
$ws.Cells.Item(7, 3).Value = -0.123754865325576
$ws.Cells.Item(7, 4).Value = -0.8027500595354817
$ws.Cells.Item(7, 5).Value = -0.4889830817910152
$ws.Cells.Item(7, 6).Value = 0.6700000762939453
$ws.Cells.Item(7, 7).Value = 1.139506816864014
$ws.Cells.Item(7, 8).Value = 2.569955825805664
$ws.Cells.Item(7, 9).Value = 1.812661170959473

$ws.Cells.Item(8, 1).Value = "model_9_2_18"
$ws.Cells.Item(8, 2).Value = 0.394807147247453
$ws.Cells.Item(8, 3).Value = -0.1209613893088219
$ws.Cells.Item(8, 4).Value = -0.8025171828478419
$ws.Cells.Item(8, 5).Value = -0.4876209277160168
$ws.Cells.Item(8, 6).Value = 0.6697701811790466
$ws.Cells.Item(8, 7).Value = 1.136674404144287
$ws.Cells.Item(8, 8).Value = 2.569623470306396
$ws.Cells.Item(8, 9).Value = 1.811002850532532

$ws.Cells.Item(9, 1).Value = "model_9_2_17"
$ws.Cells.Item(9, 2).Value = 0.4083101509329208
$ws.Cells.Item(9, 3).Value = -0.1429893640532816
$ws.Cells.Item(9, 4).Value = -0.7034792147748266
$ws.Cells.Item(9, 5).Value = -0.4427573905768438
$ws.Cells.Item(9, 6).Value = 0.654826283454895
$ws.Cells.Item(9, 7).Value = 1.159011125564575
$ws.Cells.Item(9, 8).Value = 2.42843770980835
$ws.Cells.Item(9, 9).Value = 1.756386876106262

$ws.Cells.Item(10, 1).Value = "model_9_2_15"
$ws.Cells.Item(10, 2).Value = 0.4162104894360158
$ws.Cells.Item(10, 3).Value = -0.2533067777977969
$ws.Cells.Item(10, 4).Value = -0.5415231811925316
$ws.Cells.Item(10, 5).Value = -0.4021564811833105
$ws.Cells.Item(10, 6).Value = 0.6460829377174377
$ws.Cells.Item(10, 7).Value = 1.270874857902527
$ws.Cells.Item(10, 8).Value = 2.197556972503662
$ws.Cells.Item(10, 9).Value = 1.706960082054138

$ws.Cells.Item(11, 1).Value = "model_9_2_16"
$ws.Cells.Item(11, 2).Value = 0.4171246948378288
$ws.Cells.Item(11, 3).Value = -0.1776519777552443
$ws.Cells.Item(11, 4).Value = -0.6157892313930369
$ws.Cells.Item(11, 5).Value = -0.4097162239708465
$ws.Cells.Item(11, 6).Value = 0.6450712084770203
$ws.Cells.Item(11, 7).Value = 1.194159626960754
$ws.Cells.Item(11, 8).Value = 2.303428888320923
$ws.Cells.Item(11, 9).Value = 1.716163277626038

$ws.Cells.Item(12, 1).Value = "model_9_2_14"
$ws.Cells.Item(12, 2).Value = 0.4178975388947209
$ws.Cells.Item(12, 3).Value = -0.3218548100673551
$ws.Cells.Item(12, 4).Value = -0.4567354599741673
$ws.Cells.Item(12, 5).Value = -0.3856598886957869
$ws.Cells.Item(12, 6).Value = 0.6442158818244934
$ws.Cells.Item(12, 7).Value = 1.340383648872375
$ws.Cells.Item(12, 8).Value = 2.076685905456543
$ws.Cells.Item(12, 9).Value = 1.686877489089966

$ws.Cells.Item(13, 1).Value = "model_9_2_0"
$ws.Cells.Item(13, 2).Value = 0.4313727293114354
$ws.Cells.Item(13, 3).Value = 0.5331643283309573
$ws.Cells.Item(13, 4).Value = 0.85934028293669
$ws.Cells.Item(13, 5).Value = 0.7166254150836009
$ws.Cells.Item(13, 6).Value = 0.6293028593063354
$ws.Cells.Item(13, 7).Value = 0.4733794629573822
$ws.Cells.Item(13, 8).Value = 0.2005209922790527
$ws.Cells.Item(13, 9).Value = 0.3449751436710358

$ws.Cells.Item(14, 1).Value = "model_9_2_13"
$ws.Cells.Item(14, 2).Value = 0.4425568157026716
$ws.Cells.Item(14, 3).Value = -0.4315412823431584
$ws.Cells.Item(14, 4).Value = -0.2044962957137921
$ws.Cells.Item(14, 5).Value = -0.2950249903333444
$ws.Cells.Item(14, 6).Value = 0.6169253587722778
$ws.Cells.Item(14, 7).Value = 1.451607584953308
$ws.Cells.Item(14, 8).Value = 1.717099905014038
$ws.Cells.Item(14, 9).Value = 1.576540231704712

$ws.Cells.Item(15, 1).Value = "model_9_2_12"
$ws.Cells.Item(15, 2).Value = 0.4478224952346385
$ws.Cells.Item(15, 3).Value = -0.4686402154567268
$ws.Cells.Item(15, 4).Value = -0.1064846936315609
$ws.Cells.Item(15, 5).Value = -0.2573768615676526
$ws.Cells.Item(15, 6).Value = 0.6110978126525879
$ws.Cells.Item(15, 7).Value = 1.489226698875427
$ws.Cells.Item(15, 8).Value = 1.577377080917358
$ws.Cells.Item(15, 9).Value = 1.530707955360413

$ws.Cells.Item(16, 1).Value = "model_9_2_1"
$ws.Cells.Item(16, 2).Value = 0.4612245651096805
$ws.Cells.Item(16, 3).Value = 0.4098655412951309
$ws.Cells.Item(16, 4).Value = 0.7747377342177717
$ws.Cells.Item(16, 5).Value = 0.6156318690909107
$ws.Cells.Item(16, 6).Value = 0.5962656736373901
$ws.Cells.Item(16, 7).Value = 0.5984065532684326
$ws.Cells.Item(16, 8).Value = 0.3211282789707184
$ws.Cells.Item(16, 9).Value = 0.4679228663444519

$ws.Cells.Item(17, 1).Value = "model_9_2_10"
$ws.Cells.Item(17, 2).Value = 0.4682376824657087
$ws.Cells.Item(17, 3).Value = -0.4122061538720287
$ws.Cells.Item(17, 4).Value = 0.04860288004699442
$ws.Cells.Item(17, 5).Value = -0.1470287389671678
$ws.Cells.Item(17, 6).Value = 0.5885042548179626
$ws.Cells.Item(17, 7).Value = 1.43200159072876
$ws.Cells.Item(17, 8).Value = 1.356287956237793
$ws.Cells.Item(17, 9).Value = 1.396372199058533

$ws.Cells.Item(18, 1).Value = "model_9_2_11"
$ws.Cells.Item(18, 2).Value = 0.4685733937419594
$ws.Cells.Item(18, 3).Value = -0.4150521055123091
$ws.Cells.Item(18, 4).Value = -0.001305962492374135
$ws.Cells.Item(18, 5).Value = -0.1757868094377395
$ws.Cells.Item(18, 6).Value = 0.5881325602531433
$ws.Cells.Item(18, 7).Value = 1.434887409210205
$ws.Cells.Item(18, 8).Value = 1.427436828613281
$ws.Cells.Item(18, 9).Value = 1.431381821632385

$ws.Cells.Item(19, 1).Value = "model_9_2_5"
$ws.Cells.Item(19, 2).Value = 0.4770748840010981
$ws.Cells.Item(19, 3).Value = -0.6785552249647064
$ws.Cells.Item(19, 4).Value = 0.6192785891760353
$ws.Cells.Item(19, 5).Value = 0.04999917295311584
$ws.Cells.Item(19, 6).Value = 0.5787240266799927
$ws.Cells.Item(19, 7).Value = 1.702084064483643
$ws.Cells.Item(19, 8).Value = 0.5427469611167908
$ws.Cells.Item(19, 9).Value = 1.156513929367065

$ws.Cells.Item(20, 1).Value = "model_9_2_2"
$ws.Cells.Item(20, 2).Value = 0.4777152792038644
$ws.Cells.Item(20, 3).Value = 0.4710855768393517
$ws.Cells.Item(20, 4).Value = 0.7731618305924106
$ws.Cells.Item(20, 5).Value = 0.6417600922738224
$ws.Cells.Item(20, 6).Value = 0.5780153274536133
$ws.Cells.Item(20, 7).Value = 0.5363283753395081
$ws.Cells.Item(20, 8).Value = 0.32337486743927
$ws.Cells.Item(20, 9).Value = 0.4361148476600647

$ws.Cells.Item(21, 1).Value = "model_9_2_3"
$ws.Cells.Item(21, 2).Value = 0.4859424368311136
$ws.Cells.Item(21, 3).Value = 0.3431803431843686
$ws.Cells.Item(21, 4).Value = 0.6055340807479523
$ws.Cells.Item(21, 5).Value = 0.4929828178031755
$ws.Cells.Item(21, 6).Value = 0.5689102411270142
$ws.Cells.Item(21, 7).Value = 0.6660264730453491
$ws.Cells.Item(21, 8).Value = 0.5623407959938049
$ws.Cells.Item(21, 9).Value = 0.6172336339950562

$ws.Cells.Item(22, 1).Value = "model_9_2_9"
$ws.Cells.Item(22, 2).Value = 0.494209294178147
$ws.Cells.Item(22, 3).Value = -0.3168884535899943
$ws.Cells.Item(22, 4).Value = 0.1713427398360524
$ws.Cells.Item(22, 5).Value = -0.03735762101433471
$ws.Cells.Item(22, 6).Value = 0.5597612857818604
$ws.Cells.Item(22, 7).Value = 1.335347652435303
$ws.Cells.Item(22, 8).Value = 1.181313157081604
$ws.Cells.Item(22, 9).Value = 1.262860536575317

$ws.Cells.Item(23, 1).Value = "model_9_2_8"
$ws.Cells.Item(23, 2).Value = 0.5125015733478058
$ws.Cells.Item(23, 3).Value = -0.2364466407462877
$ws.Cells.Item(23, 4).Value = 0.217294714807736
$ws.Cells.Item(23, 5).Value = 0.02343745320494239
$ws.Cells.Item(23, 6).Value = 0.5395171046257019
$ws.Cells.Item(23, 7).Value = 1.253778219223022
$ws.Cells.Item(23, 8).Value = 1.115805149078369
$ws.Cells.Item(23, 9).Value = 1.188849687576294

$ws.Cells.Item(24, 1).Value = "model_9_2_4"
$ws.Cells.Item(24, 2).Value = 0.5144264659964439
$ws.Cells.Item(24, 3).Value = -0.3853943535076529
$ws.Cells.Item(24, 4).Value = 0.6677973452548888
$ws.Cells.Item(24, 5).Value = 0.2060127544917079
$ws.Cells.Item(24, 6).Value = 0.5373868346214294
$ws.Cells.Item(24, 7).Value = 1.404813885688782
$ws.Cells.Item(24, 8).Value = 0.4735798239707947
$ws.Cells.Item(24, 9).Value = 0.9665859341621399

$ws.Cells.Item(25, 1).Value = "model_9_2_6"
$ws.Cells.Item(25, 2).Value = 0.5187547868627813
$ws.Cells.Item(25, 3).Value = -0.212201553680553
$ws.Cells.Item(25, 4).Value = 0.4462686473691653
$ws.Cells.Item(25, 5).Value = 0.1603086173375706
$ws.Cells.Item(25, 6).Value = 0.5325966477394104
$ws.Cells.Item(25, 7).Value = 1.229193329811096
$ws.Cells.Item(25, 8).Value = 0.7893856167793274
$ws.Cells.Item(25, 9).Value = 1.022225260734558

$ws.Cells.Item(26, 1).Value = "model_9_2_7"
$ws.Cells.Item(26, 2).Value = 0.5387375338087419
$ws.Cells.Item(26, 3).Value = -0.1528337384252552
$ws.Cells.Item(26, 4).Value = 0.3989360722626639
$ws.Cells.Item(26, 5).Value = 0.1604052180980393
$ws.Cells.Item(26, 6).Value = 0.5104816555976868
$ws.Cells.Item(26, 7).Value = 1.168993473052979
$ws.Cells.Item(26, 8).Value = 0.8568617701530457
$ws.Cells.Item(26, 9).Value = 1.022107601165771
